# Add FTP connection functionality: duplicate the config sheet into a new
# "FUSE" sheet (placed after the existing sheet) and make it the active tab,
# with its own selection at L23. The original sheet's selection (C23) is
# left untouched and it is no longer the selected tab.

$wb = $excel.ActiveWorkbook
$sourceSheet = $wb.Worksheets.Item(1)

# Copy the existing sheet, placing the new copy right after it.
$sourceSheet.Copy($null, $sourceSheet)

# The newly created copy becomes the last sheet in the workbook.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "FUSE"

# Give the new sheet its own selected cell / active tab state.
$newSheet.Range("L23").Select()
